$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Lower Right Cell" column (D) for the scenario-indices block (rows 5-11)
# previously pointed at row 17 of the referenced ranges; the variable block
# now starts one row later, so these need to point at row 18 instead.
$ws.Range("D5").Value  = "A18"
$ws.Range("D6").Value  = "B18"
$ws.Range("D7").Value  = "C18"
$ws.Range("D8").Value  = "G18"
$ws.Range("D9").Value  = "H18"
$ws.Range("D10").Value = "I18"
$ws.Range("D11").Value = "J18"

# Move the active selection from D16 to D12
$ws.Range("D12").Select() | Out-Null
